$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

# C1 / D1: drop the inherited bold/box/centered style and give them fresh
# "plain font, thin top+bottom" (C1) / "plain font, thin top+bottom+right" (D1)
# borders - this reproduces the two new cellXfs entries (borderId 4 and 5).
$ws1.Range("C1").ClearFormats()
$ws1.Range("C1").Borders.Item(8).LineStyle = 1
$ws1.Range("C1").Borders.Item(9).LineStyle = 1

$ws1.Range("D1").ClearFormats()
$ws1.Range("D1").Borders.Item(8).LineStyle = 1
$ws1.Range("D1").Borders.Item(9).LineStyle = 1
$ws1.Range("D1").Borders.Item(10).LineStyle = 1

# Anonymize the "fedcore" label
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

$ws2.Range("C1").ClearFormats()
$ws2.Range("C1").Borders.Item(8).LineStyle = 1
$ws2.Range("C1").Borders.Item(9).LineStyle = 1

$ws2.Range("D1").ClearFormats()
$ws2.Range("D1").Borders.Item(8).LineStyle = 1
$ws2.Range("D1").Borders.Item(9).LineStyle = 1
$ws2.Range("D1").Borders.Item(10).LineStyle = 1

$ws2.Range("F1").ClearFormats()
$ws2.Range("F1").Borders.Item(8).LineStyle = 1
$ws2.Range("F1").Borders.Item(9).LineStyle = 1

$ws2.Range("G1").ClearFormats()
$ws2.Range("G1").Borders.Item(8).LineStyle = 1
$ws2.Range("G1").Borders.Item(9).LineStyle = 1
$ws2.Range("G1").Borders.Item(10).LineStyle = 1

# Anonymize the "fedcore" labels
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
